$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column A content paths (rows 2-11)
$ws.Range("A2").Value = "/content/abbvie-pro/uk/en/test-podcast/podcast-examples"
$ws.Range("A3").Value = "/content/abbvie-pro/uk/en/test-podcast/tabs-test-page"
$ws.Range("A4").Value = "/content/abbvie-pro/uk/en/test-podcast/test-gif"
$ws.Range("A5").Value = "/content/abbvie-pro/uk/en/test-podcast/test-logo-lozenge-combination"
$ws.Range("A6").Value = "/content/abbvie-pro/uk/en/neuroscience/parkinsons/products/duodopa-what-is-apd/contact-us/abc"
$ws.Range("A7").Value = "/content/abbvie-pro/uk/en/test-podcast/test-aml-page"
$ws.Range("A8").Value = "/content/abbvie-pro/uk/en/test-podcast/understanding-ndo"
$ws.Range("A9").Value = "/content/abbvie-pro/uk/en/neuroscience/migraine/aquipta/contact-us1/error"
$ws.Range("A10").Value = "/content/abbvie-pro/uk/en/test-podcast/dosing-calculator"
$ws.Range("A11").Value = "/content/abbvie-pro/uk/en/test-podcast/test-header-and-footer-addition"

# Clear the Status cells that no longer have values
$ws.Range("B2").ClearContents()
$ws.Range("B6").ClearContents()

# Set the Status cells that are now populated
$ws.Range("B7").Value = "Preview activated"
$ws.Range("B11").Value = "Preview activated"

# Add new rows 12-14
$ws.Range("A12").Value = "/content/abbvie-pro/uk/en/test-podcast/understanding-oab"
$ws.Range("B12").Value = "Preview activated"

$ws.Range("A13").Value = "/content/abbvie-pro/uk/en/test-podcast/cll-life-test"
$ws.Range("B13").Value = "Preview activated"

$ws.Range("A14").Value = "/content/abbvie-pro/uk/en/test-podcast/rinvoq-gastro-home-page-template"
$ws.Range("B14").Value = "Preview activated"
